# Weekly crime data refresh (new week of data collected).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / masthead text updates -----------------------------------
# Mayor's name (shared string used by cell M6)
$ws.Range("M6").Value = "Thomas G. Donlon"

# "Volume 31   Number  38" -> "... Number  39"
$ws.Range("A8").Value = "Volume 31   Number  39"

# "Report Covering the Week  9/16/2024  Through  9/22/2024" -> next week
$ws.Range("C9").Value = "Report Covering the Week  9/23/2024  Through  9/29/2024"

# --- Cells that change between a numeric value and a text marker ------
# ("0" and "***.*" are stored as literal text in this workbook, not as
# numbers, so we clone the style+content of an untouched donor cell that
# already carries the right text, then overwrite numeric cells as needed.)

# D15 : 1 -> "0"
$ws.Range("C15").Copy($ws.Range("D15"))
# E15 : -100 -> "***.*"
$ws.Range("N22").Copy($ws.Range("E15"))
# C17 : 3 -> "0"
$ws.Range("C15").Copy($ws.Range("C17"))
# C22 : "0" -> 1 (numeric)
$ws.Range("F15").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
# D22 : 1 -> "0"
$ws.Range("C15").Copy($ws.Range("D22"))
# E22 : -100 -> "***.*"
$ws.Range("N22").Copy($ws.Range("E22"))
# D27 : 1 -> "0"
$ws.Range("C15").Copy($ws.Range("D27"))
# E27 : -100 -> "***.*"
$ws.Range("N22").Copy($ws.Range("E27"))
# F29 : 1 -> "0"
$ws.Range("C15").Copy($ws.Range("F29"))
# F30 : 1 -> "0"
$ws.Range("C15").Copy($ws.Range("F30"))

# --- Plain numeric value updates across the crime-stat table ----------
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 79
$ws.Range("J16").Value = 102
$ws.Range("K16").Value = -22.549019607843
$ws.Range("L16").Value = -42.753623188405
$ws.Range("M16").Value = -24.038461538461
$ws.Range("N16").Value = -86.877076411960

$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = -41.176470588235
$ws.Range("J17").Value = 165
$ws.Range("K17").Value = 3.636363636363
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 69.306930693069
$ws.Range("N17").Value = 0.588235294117

$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -73.333333333333
$ws.Range("I18").Value = 91
$ws.Range("J18").Value = 132
$ws.Range("K18").Value = -31.060606060606
$ws.Range("L18").Value = -32.089552238806
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = -68.181818181818

$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 15.384615384615
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -3.703703703703
$ws.Range("I19").Value = 375
$ws.Range("J19").Value = 444
$ws.Range("K19").Value = -15.540540540540
$ws.Range("L19").Value = -36.224489795918
$ws.Range("M19").Value = 88.442211055276
$ws.Range("N19").Value = 9.011627906976

$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 42
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = -4.545454545454
$ws.Range("L20").Value = -19.230769230769
$ws.Range("M20").Value = 2.439024390243
$ws.Range("N20").Value = -85.467128027681

$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -20.833333333333
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -21.359223300970
$ws.Range("I21").Value = 767
$ws.Range("J21").Value = 897
$ws.Range("K21").Value = -14.492753623188
$ws.Range("L21").Value = -30.145719489981
$ws.Range("M21").Value = 46.653919694072
$ws.Range("N21").Value = -55.250875145857

$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 18
$ws.Range("K22").Value = 28.571428571428
$ws.Range("L22").Value = -10
$ws.Range("M22").Value = 38.461538461538

$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -83.333333333333
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = -31.578947368421
$ws.Range("I23").Value = 155
$ws.Range("J23").Value = 126
$ws.Range("K23").Value = 23.015873015873
$ws.Range("L23").Value = 11.510791366906
$ws.Range("M23").Value = 56.565656565656

$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 93.333333333333
$ws.Range("F24").Value = 120
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 46.341463414634
$ws.Range("I24").Value = 1033
$ws.Range("J24").Value = 925
$ws.Range("K24").Value = 11.675675675675
$ws.Range("L24").Value = -42.096412556053
$ws.Range("M24").Value = 85.791366906474

$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 214.285714285714
$ws.Range("F25").Value = 85
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 102.380952380952
$ws.Range("I25").Value = 696
$ws.Range("J25").Value = 510
$ws.Range("K25").Value = 36.470588235294
$ws.Range("L25").Value = -51.362683438155

$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -7.692307692307
$ws.Range("I26").Value = 322
$ws.Range("J26").Value = 337
$ws.Range("K26").Value = -4.451038575667
$ws.Range("L26").Value = -6.122448979591
$ws.Range("M26").Value = 28.8

$ws.Range("G28").Value = 4
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -16.666666666666

$ws.Range("L29").Value = -60

$ws.Range("L30").Value = -50

$ws.Range("I33").Value = 6
$ws.Range("L33").Value = 500
